# Edit /tmp/work/before.pptx to match the target commit:
#   - Slide 1, "Rectangle 4" text box: second paragraph "Test Report" ->
#     "Editing this line on tarakeena account" (typed as three runs, the
#     way PowerPoint's editor splits a run when new text is appended
#     after an in-place edit).
#   - Slide 1, "Table 1": Version cell "1.1" -> "1.2".

$p = $ppt.ActivePresentation
$s = $p.Slides.Item(1)

# --- 1) Rectangle 4 text box: "Test Report" -> "Editing this line on tarakeena account"
$shp = $s.Shapes.Item(1)
$tf = $shp.TextFrame
$tr = $tf.TextRange

$para = $tr.Paragraphs(2)
$run = $para.Runs(1)

# Re-use the existing run for the first chunk of text, then append the
# remaining text as new runs (InsertAfter keeps the same character
# formatting, matching the rPr of the original run).
$run.Text = "Editing this line on "
$run.InsertAfter("tarakeena") | Out-Null

$para = $tr.Paragraphs(2)
$lastRun = $para.Runs($para.Runs().Count)
$lastRun.InsertAfter(" account") | Out-Null

# --- 2) Table 1: Version cell "1.1" -> "1.2"
$tblShape = $s.Shapes.Item(2)
$tbl = $tblShape.Table
$cell = $tbl.Cell(4, 2)
$cell.Shape.TextFrame.TextRange.Runs(1).Text = "1.2"
